$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.199.63'
$ws.Range('E2').Value = '  +0.68%  '
$ws.Range('D3').Value = '2.095.81'
$ws.Range('E3').Value = '  +2.97%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''229.31'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('E6').Value = '  +0.55%  '
$ws.Range('D7').Value = '''60.85'
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').Value = '''0.0844'
$ws.Range('E10').Value = '  +2.54%  '
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('D12').Value = '2.405.45'
$ws.Range('E12').Value = '  +2.98%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = '''14.67'
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '''22.32'
$ws.Range('E14').Value = '  +4.25%  '
$ws.Range('D15').Value = '''5.50'
$ws.Range('E15').Value = '  +6.48%  '
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('D17').Value = '2.078.47'
$ws.Range('E17').Value = '  +1.29%  '
$ws.Range('D18').Value = '38.145.12'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('E19').Value = '  +1.84%  '
$ws.Range('D20').Value = '''70.23'
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('E21').Value = '  +1.01%  '
$ws.Range('D22').Value = '''223.87'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('E24').Value = '  -0.92%  '
$ws.Range('E25').Value = '  +3.62%  '
$ws.Range('D26').Value = '''170.27'
$ws.Range('E26').Value = '  +2.03%  '
$ws.Range('D27').Value = '''9.46'
$ws.Range('E27').Value = '  +1.70%  '
$ws.Range('E28').Value = '  +0.61%  '
$ws.Range('D29').Value = '''19.01'
$ws.Range('E29').Value = '  +0.83%  '
$ws.Range('E30').Value = '  +6.34%  '
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('E32').Value = '  +4.41%  '
$ws.Range('E33').Value = '  +4.04%  '
$ws.Range('D34').Value = '''4.43'
$ws.Range('E34').Value = '  +0.56%  '
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('D36').Value = '''6.47'
$ws.Range('E36').Value = '  +1.64%  '
$ws.Range('E37').Value = '  +5.16%  '
$ws.Range('D38').Value = '''3.56'
$ws.Range('E38').Value = '  +8.53%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').Value = '''18.03'
$ws.Range('E40').Value = '  +2.07%  '
$ws.Range('D41').Value = '1.559.01'
$ws.Range('E41').Value = '  +1.27%  '
$ws.Range('D42').Value = '''100.41'
$ws.Range('E42').Value = '  +4.30%  '
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('E44').Value = '  +1.48%  '
$ws.Range('E45').Value = '  -0.63%  '
$ws.Range('D46').Value = '''4.15'
$ws.Range('E46').Value = '  +3.63%  '
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('E48').Value = '  +1.72%  '
$ws.Range('D49').Value = '''7.27'
$ws.Range('E49').Value = '  +2.76%  '
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('D51').Value = '2.292.82'
$ws.Range('E51').Value = '  +3.02%  '
